$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: orderedOn, name, email, phone, products, quantity, colour, totalAmount, paymentMethod, orderedStatus
$data = @(
    @("Wed Feb 15 2023 ", "trideep kumar", "trideep@gmail.com", "9061117489", "Casio G-Shock Gold Dial Mens Watch-G1053", 4, "Black and Gold", 76030, "RAZORPAY", "Placed"),
    @("Fri Feb 17 2023 ", "Mazin ", "mazinshajahan4444@gmail.com", "7510722928", "Casio G-Shock Gold Dial Mens Watch-G1053", 1, "Black and Gold", 19045, "COD", "Delivered"),
    @("Sat Feb 18 2023 ", "Muhammed Musthafa", "musthafa723262@gmail.com", "9544535049", "Mens 45.9 mm Daniel Black Dial Zinc Alloy Watch - NCTH1710383", 1, "black", 13300, "Wallet", "Placed"),
    @("Sat Feb 18 2023 ", "Muhammed Musthafa", "musthafa723262@gmail.com", "9544535049", "Mens 45.9 mm Daniel Black Dial Zinc Alloy Watch - NCTH1710383", 1, "black", 13300, "Wallet", "Delivered"),
    @("Sat Feb 18 2023 ", "Muhammed Musthafa", "musthafa723262@gmail.com", "9544535049", "Casio Edifice Black Dial Men's Watch -EX511", 1, "silver", 30045, "Wallet", "Delivered"),
    @("Sat Feb 18 2023 ", "Muhammed Musthafa", "musthafa723262@gmail.com", "9544535049", "Mens 45.9 mm Daniel Black Dial Zinc Alloy Watch - NCTH1710383", 1, "black", 13300, "Wallet", "Placed"),
    @("Sat Feb 18 2023 ", "Muhammed Musthafa", "musthafa723262@gmail.com", "9544535049", "Mens 45.9 mm Daniel Black Dial Zinc Alloy Watch - NCTH1710383", 1, "black", 13300, "Wallet", "Placed"),
    @("Sat Feb 18 2023 ", "Muhammed Musthafa", "musthafa723262@gmail.com", "9544535049", "Mens 45.9 mm Daniel Black Dial Zinc Alloy Watch - NCTH1710383", 2, "black", 20455, "Wallet", "Placed"),
    @("Sat Feb 18 2023 ", "Muhammed Musthafa", "musthafa723262@gmail.com", "9544535049", "Mens 45.9 mm Daniel Black Dial Zinc Alloy Watch - NCTH1710383", 1, "black", 13300, "Wallet", "Placed"),
    @("Sat Feb 18 2023 ", "Muhammed Musthafa", "musthafa723262@gmail.com", "9544535049", "Mens 45.9 mm Daniel Black Dial Zinc Alloy Watch - NCTH1710383", 1, "black", 11710, "Wallet", "Placed"),
    @("Sat Feb 18 2023 ", "Muhammed Musthafa", "musthafa723262@gmail.com", "9544535049", "Mens 45.9 mm Daniel Black Dial Zinc Alloy Watch - NCTH1710383", 1, "black", 11710, "Wallet", "Placed"),
    @("Sat Feb 18 2023 ", "Musthafa", "musthafa916@gmail.com", "9047792111", "DAY-DATE 40", 2, " yellow gold and diamonds", 3136350, "COD", "Delivered"),
    @("Mon Feb 20 2023 ", "neeraj", "ifoxscand@gmail.com", "9074165714", "Casio G-Shock Gold Dial Mens Watch-G1053", 1, "Black and Gold", 19045, "COD", "Delivered"),
    @("Mon Feb 20 2023 ", "neeraj", "ifoxscand@gmail.com", "9074165714", "Mens 45.9 mm Daniel Black Dial Zinc Alloy Watch - NCTH1710383", 1, "black", 11710, "RAZORPAY", "Placed")
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    # phone numbers must remain text, not get coerced into numbers
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[3]

    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
}
